$d = $word.ActiveDocument
$q1 = [char]0x201C
$q2 = [char]0x201D

# -----------------------------------------------------------------
# Append six new paragraphs after the final paragraph of the body
# (the one ending in "...six phases of the design process." which
# also carries the _GoBack bookmark). The bookmark stays attached to
# that paragraph; new paragraphs are appended after it.
#
# Note: InsertParagraphAfter() copies the pPr (including pStyle) of
# the paragraph it is called on into the freshly created paragraph,
# so every new paragraph's style is set explicitly right after it is
# created (Normal is set explicitly too, which the engine renders
# without an explicit <w:pStyle> tag, matching default Normal
# paragraphs elsewhere in the document).
# -----------------------------------------------------------------

# 1) Heading2: "Multiple categories"
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Style = $d.Styles.Item("Heading 2")
$pr = $d.Paragraphs.Item($idx).Range
$pr.Text = "Multiple categories"

# 2) Normal: "After switching to design process questions..."
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Style = $d.Styles.Item("Normal")
$pr = $d.Paragraphs.Item($idx).Range
$pr.Text = "After switching to design process questions, I added a category for each step in the design process. The player can take as many or as few categories as they choose before ending the quiz and completing the game."

# 3) Heading1: "Making"
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Style = $d.Styles.Item("Heading 1")
$pr = $d.Paragraphs.Item($idx).Range
$pr.Text = "Making"

# 4) Normal: "With the engine completed..."
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Style = $d.Styles.Item("Normal")
$pr = $d.Paragraphs.Item($idx).Range
$pr.Text = "With the engine completed, I now had to program in all of the questions. I wanted to have at least 5 questions in each category. With 6 categories this meant 30 questions in total. This was programmed in using the existing engine. Additionally, graphical effects were introduced such as dashes above and below text, and display of which category each question was a part of."

# 5) Heading1: "Sharing"
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Style = $d.Styles.Item("Heading 1")
$pr = $d.Paragraphs.Item($idx).Range
$pr.Text = "Sharing"

# 6) Normal: two runs of feedback text
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Style = $d.Styles.Item("Normal")
$pr = $d.Paragraphs.Item($idx).Range
$pr.Text = "After sharing the project with the other Programming Level 2 students, the key feedback was the limitation of the programs input. The program only allowed the user to input capital letters. While I would have liked to fix this issue, "
$pr2 = $d.Paragraphs.Item($idx).Range
$pr2.Collapse(0)
$pr2.InsertAfter("it is rooted in the main game engine using a string to represent the correct answer ($($q1)A$($q2)). To fix the issue, a conversion would need to be made from uppercase to lowercase letters which is not easily documented. I will begin work to make the improvement.")
